$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.204854011535645
$ws.Range("B1").Value = 2.123377561569214
$ws.Range("C1").Value = 5.835978031158447
$ws.Range("D1").Value = 1.043137431144714
$ws.Range("E1").Value = 1.199688792228699
